$wb = $excel.ActiveWorkbook

# Remove the "Sheet1" worksheet (fruit pricing scratch sheet no longer needed).
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet1").Delete()

# Rename the remaining sheet.
$ws = $wb.Worksheets.Item("payment-request")
$ws.Name = "testSheet"
